# Generate Report for Handoff
# Rewrites the Overview / zh-cn / de-de sheets so that the two new files
# (3e733c7f-...png and bc2883ff-...png, both dependencies of the renamed
# a91dee86-...md source file) show up alongside the existing
# .localization-config row, and refreshes the handoff timestamps.

$wb = $excel.ActiveWorkbook

$DATE_FMT = "yyyy-mm-dd HH:mm:ss"

# NOTE: cell styling for the linked cells is applied later, in one pass,
# via Hyperlinks.Add (which stamps its own HyperLink look onto the cell) -
# setting Font properties here would just get clobbered by that call.
function Set-LinkCell($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
}

function Set-DateCell($ws, $addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = $DATE_FMT
    $r.Value = $text
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("A1").Value = "File Name"
$ws1.Range("B1").Value = "zh-cn"
$ws1.Range("C1").Value = "de-de"

Set-LinkCell $ws1 "A2" "3e733c7f-fb97-4e3d-84a9-46acce90192e.png"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

Set-LinkCell $ws1 "A3" "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

Set-LinkCell $ws1 "A4" "bc2883ff-a2da-460e-8416-54d43d14eaab.png"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

Set-LinkCell $ws1 "A5" ".localization-config"
$ws1.Range("B5").Value = "Not to be localized"
$ws1.Range("C5").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/3e733c7f-fb97-4e3d-84a9-46acce90192e.png", "", "", "3e733c7f-fb97-4e3d-84a9-46acce90192e.png") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md", "", "", "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/bc2883ff-a2da-460e-8416-54d43d14eaab.png", "", "", "bc2883ff-a2da-460e-8416-54d43d14eaab.png") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("A1").Value = "Source File Name"
$ws2.Range("B1").Value = "Status"
$ws2.Range("C1").Value = "Latest Handoff File"
$ws2.Range("D1").Value = "Latest Handoff Datetime"
$ws2.Range("E1").Value = "Latest Target File"
$ws2.Range("F1").Value = "Latest Handback File"
$ws2.Range("G1").Value = "Latest Handback DateTime"
$ws2.Range("H1").Value = "Handoff Reason"
$ws2.Range("I1").Value = "Dependency From"

# Row 2: new dependency file (png)
Set-LinkCell $ws2 "A2" "3e733c7f-fb97-4e3d-84a9-46acce90192e.png"
$ws2.Range("B2").Value = "Ready for handoff"
Set-LinkCell $ws2 "C2" "0eb7870a69c5a3d915373054ec32eac193631a41.png"
Set-DateCell $ws2 "D2" "2016-03-10 03:42:37"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "IsDependency"
$ws2.Range("I2").Value = "e2e\a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md"

# Row 3: renamed source markdown file
Set-LinkCell $ws2 "A3" "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md"
$ws2.Range("B3").Value = "Ready for handoff"
Set-LinkCell $ws2 "C3" "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.cad963e9905e2cfea3814a8405f93e9d5f791a61.zh-cn.xlf"
Set-DateCell $ws2 "D3" "2016-03-10 03:42:37"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

# Row 4: new dependency file (png)
Set-LinkCell $ws2 "A4" "bc2883ff-a2da-460e-8416-54d43d14eaab.png"
$ws2.Range("B4").Value = "Ready for handoff"
Set-LinkCell $ws2 "C4" "f616413dff660d14bc66ab554864bb6ff73ccaa2.png"
Set-DateCell $ws2 "D4" "2016-03-10 03:42:37"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "IsDependency"
$ws2.Range("I4").Value = "e2e\a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md"

# Row 5: unchanged .localization-config
Set-LinkCell $ws2 "A5" ".localization-config"
$ws2.Range("B5").Value = "Not to be localized"
Set-DateCell $ws2 "D5" "0001-01-01 00:00:00"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/3e733c7f-fb97-4e3d-84a9-46acce90192e.png", "", "", "3e733c7f-fb97-4e3d-84a9-46acce90192e.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/12486af9ff48cbf2793147cdcd9bdf806d9b431d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0eb7870a69c5a3d915373054ec32eac193631a41.png", "", "", "0eb7870a69c5a3d915373054ec32eac193631a41.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md", "", "", "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/12486af9ff48cbf2793147cdcd9bdf806d9b431d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a91dee86-1d2c-4a20-bf6e-f7995814bcf1.cad963e9905e2cfea3814a8405f93e9d5f791a61.zh-cn.xlf", "", "", "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.cad963e9905e2cfea3814a8405f93e9d5f791a61.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/bc2883ff-a2da-460e-8416-54d43d14eaab.png", "", "", "bc2883ff-a2da-460e-8416-54d43d14eaab.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/12486af9ff48cbf2793147cdcd9bdf806d9b431d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f616413dff660d14bc66ab554864bb6ff73ccaa2.png", "", "", "f616413dff660d14bc66ab554864bb6ff73ccaa2.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("A1").Value = "Source File Name"
$ws3.Range("B1").Value = "Status"
$ws3.Range("C1").Value = "Latest Handoff File"
$ws3.Range("D1").Value = "Latest Handoff Datetime"
$ws3.Range("E1").Value = "Latest Target File"
$ws3.Range("F1").Value = "Latest Handback File"
$ws3.Range("G1").Value = "Latest Handback DateTime"
$ws3.Range("H1").Value = "Handoff Reason"
$ws3.Range("I1").Value = "Dependency From"

# Row 2: new dependency file (png)
Set-LinkCell $ws3 "A2" "3e733c7f-fb97-4e3d-84a9-46acce90192e.png"
$ws3.Range("B2").Value = "Ready for handoff"
Set-LinkCell $ws3 "C2" "0eb7870a69c5a3d915373054ec32eac193631a41.png"
Set-DateCell $ws3 "D2" "2016-03-10 03:42:41"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "IsDependency"
$ws3.Range("I2").Value = "e2e\a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md"

# Row 3: renamed source markdown file
Set-LinkCell $ws3 "A3" "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md"
$ws3.Range("B3").Value = "Ready for handoff"
Set-LinkCell $ws3 "C3" "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.cad963e9905e2cfea3814a8405f93e9d5f791a61.de-de.xlf"
Set-DateCell $ws3 "D3" "2016-03-10 03:42:41"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

# Row 4: new dependency file (png)
Set-LinkCell $ws3 "A4" "bc2883ff-a2da-460e-8416-54d43d14eaab.png"
$ws3.Range("B4").Value = "Ready for handoff"
Set-LinkCell $ws3 "C4" "f616413dff660d14bc66ab554864bb6ff73ccaa2.png"
Set-DateCell $ws3 "D4" "2016-03-10 03:42:41"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "IsDependency"
$ws3.Range("I4").Value = "e2e\a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md"

# Row 5: unchanged .localization-config
Set-LinkCell $ws3 "A5" ".localization-config"
$ws3.Range("B5").Value = "Not to be localized"
Set-DateCell $ws3 "D5" "0001-01-01 00:00:00"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/3e733c7f-fb97-4e3d-84a9-46acce90192e.png", "", "", "3e733c7f-fb97-4e3d-84a9-46acce90192e.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/918e6e3df7dcc1e41b6cdc79af165d00c2b7d000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0eb7870a69c5a3d915373054ec32eac193631a41.png", "", "", "0eb7870a69c5a3d915373054ec32eac193631a41.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md", "", "", "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/918e6e3df7dcc1e41b6cdc79af165d00c2b7d000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a91dee86-1d2c-4a20-bf6e-f7995814bcf1.cad963e9905e2cfea3814a8405f93e9d5f791a61.de-de.xlf", "", "", "a91dee86-1d2c-4a20-bf6e-f7995814bcf1.cad963e9905e2cfea3814a8405f93e9d5f791a61.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/e2e/bc2883ff-a2da-460e-8416-54d43d14eaab.png", "", "", "bc2883ff-a2da-460e-8416-54d43d14eaab.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/918e6e3df7dcc1e41b6cdc79af165d00c2b7d000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f616413dff660d14bc66ab554864bb6ff73ccaa2.png", "", "", "f616413dff660d14bc66ab554864bb6ff73ccaa2.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/8be002573d582ef5ca28aef650506e8aea942afb/.localization-config", "", "", ".localization-config") | Out-Null

$ws1.Range("A1").Select() | Out-Null
